$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.196.22'
$ws.Range('E2').Value = '  -0.13%  '

$ws.Range('D3').Value = '1.854.81'
$ws.Range('E3').Value = '  -0.33%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.05'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.51%  '

$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6991'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07718'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.27%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3081'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.05%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.75'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.30%  '

$ws.Range('E11').Value = '  -2.56%  '

$ws.Range('D12').Value = '1.862.66'
$ws.Range('E12').Value = '  -0.75%  '

$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '92.15'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.45%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.094'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.58%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6862'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.32%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.501'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.37%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008376'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.16%  '

$ws.Range('D18').Value = '29.232.55'
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.19'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.32%  '

$ws.Range('D20').Value = '2.115.91'
$ws.Range('E20').Value = '  -0.68%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.82'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.09%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.509'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1508'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.11%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.53'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.66%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.832'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.75%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.47'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.31%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.559'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.224'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.181'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.04%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.196'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05197'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.47%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7616'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.48%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.841'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.55%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.162'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.57%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.25%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01861'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.10%  '

$ws.Range('D39').Value = '1.214.44'
$ws.Range('E39').Value = '  -2.75%  '

$ws.Range('E40').Value = '  -0.51%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8960'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.67%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.77'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.25%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9990'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.534'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -11.86%  '

$ws.Range('D45').Value = '2.013.14'
$ws.Range('E45').Value = '  -2.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.89'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -9.52%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5182'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.38%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.516'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.35%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000121'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -5.99%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.750'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.84%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.009'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.39%  '
